$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2268.4
$ws.Range("I43").Value = 2686.75
$ws.Range("J43").Value = 595
$ws.Range("K43").Value = 2686.75
$ws.Range("L43").Value = 595
$ws.Range("M43").Value = -2617.75
$ws.Range("N43").Value = -733
$ws.Range("H107").Value = 703.7692
$ws.Range("I107").Value = 703.7692
$ws.Range("K107").Value = 703.7692
$ws.Range("M107").Value = 1216.2308
$ws.Range("H111").Value = 849
$ws.Range("J111").Value = 849
$ws.Range("L111").Value = 2547
$ws.Range("N111").Value = -8681
$ws.Range("H125").Value = 2697.4
$ws.Range("I125").Value = 2399.5
$ws.Range("K125").Value = 21595.5
$ws.Range("M125").Value = -19135.5
$ws.Range("H131").Value = 4386.0625
$ws.Range("I131").Value = 1835.2727
$ws.Range("K131").Value = 5505.8181
$ws.Range("M131").Value = -465.8181000000004
$ws.Range("H132").Value = 100019390
$ws.Range("I132").Value = 142869550
$ws.Range("K132").Value = 428608650
$ws.Range("M132").Value = -428606120
$ws.Range("H137").Value = 30305754
$ws.Range("I137").Value = 83334584
$ws.Range("K137").Value = 250003752
$ws.Range("M137").Value = -250001202
$ws.Range("H138").Value = 7834.1396
$ws.Range("I138").Value = 6636.6
$ws.Range("J138").Value = 8475.679
$ws.Range("K138").Value = 19909.8
$ws.Range("L138").Value = 25427.037
$ws.Range("M138").Value = -14769.8
$ws.Range("N138").Value = -35707.037

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 423.1111
$ws.Range("I5").Value = 423.1111
$ws.Range("K5").Value = 423.1111
$ws.Range("M5").Value = -311.1111
$ws.Range("H32").Value = 14239.223
$ws.Range("I32").Value = 14239.223
$ws.Range("K32").Value = 14239.223
$ws.Range("M32").Value = -13952.223
$ws.Range("H74").Value = 97185096
$ws.Range("J74").Value = 4375.25
$ws.Range("L74").Value = 4375.25
$ws.Range("N74").Value = -6123.25
$ws.Range("H76").Value = 53750
$ws.Range("J76").Value = 53750
$ws.Range("L76").Value = 53750
$ws.Range("N76").Value = -54426
$ws.Range("H77").Value = 97185096
$ws.Range("J77").Value = 4375.25
$ws.Range("L77").Value = 21876.25
$ws.Range("N77").Value = -30612.25
$ws.Range("H79").Value = 53750
$ws.Range("J79").Value = 53750
$ws.Range("L79").Value = 53750
$ws.Range("N79").Value = -56090
$ws.Range("H102").Value = 3288.875
$ws.Range("I102").Value = 1550.6
$ws.Range("J102").Value = 6186
$ws.Range("K102").Value = 1550.6
$ws.Range("L102").Value = 6186
$ws.Range("M102").Value = 71.40000000000009
$ws.Range("N102").Value = -9430
$ws.Range("H122").Value = 2620.875
$ws.Range("I122").Value = 2029.5
$ws.Range("J122").Value = 4395
$ws.Range("K122").Value = 6088.5
$ws.Range("L122").Value = 13185
$ws.Range("M122").Value = -3638.5
$ws.Range("N122").Value = -18085

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 423.1111
$ws.Range("I4").Value = 423.1111
$ws.Range("K4").Value = 423.1111
$ws.Range("M4").Value = -308.1111
$ws.Range("H20").Value = 2130.6155
$ws.Range("I20").Value = 1559.8
$ws.Range("K20").Value = 1559.8
$ws.Range("M20").Value = -1312.8
$ws.Range("H107").Value = 1464.7693
$ws.Range("I107").Value = 1304.9
$ws.Range("J107").Value = 1997.6666
$ws.Range("K107").Value = 1304.9
$ws.Range("L107").Value = 1997.6666
$ws.Range("M107").Value = 615.0999999999999
$ws.Range("N107").Value = -5837.6666
$ws.Range("H134").Value = 4299.6
$ws.Range("I134").Value = 2666.6667
$ws.Range("K134").Value = 8000.000100000001
$ws.Range("M134").Value = -5465.000100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2933.1333
$ws.Range("I31").Value = 2950
$ws.Range("J31").Value = 2865.6667
$ws.Range("K31").Value = 2950
$ws.Range("L31").Value = 2865.6667
$ws.Range("M31").Value = -2655
$ws.Range("N31").Value = -3455.6667
$ws.Range("H34").Value = 2933.1333
$ws.Range("I34").Value = 2950
$ws.Range("J34").Value = 2865.6667
$ws.Range("K34").Value = 2950
$ws.Range("L34").Value = 2865.6667
$ws.Range("M34").Value = -2748
$ws.Range("N34").Value = -3269.6667
$ws.Range("H99").Value = 10253.071
$ws.Range("I99").Value = 6764.5
$ws.Range("K99").Value = 6764.5
$ws.Range("M99").Value = -5266.5
$ws.Range("H126").Value = 10253.071
$ws.Range("I126").Value = 6764.5
$ws.Range("K126").Value = 20293.5
$ws.Range("M126").Value = -17823.5
$ws.Range("H141").Value = 1074824
$ws.Range("J141").Value = 1074824
$ws.Range("L141").Value = 1074824
$ws.Range("N141").Value = -1085184

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H58").Value = 37500
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H80").Value = 6983.364
$ws.Range("I80").Value = 2814
$ws.Range("J80").Value = 8546.875
$ws.Range("K80").Value = 2814
$ws.Range("L80").Value = 8546.875
$ws.Range("M80").Value = -1816
$ws.Range("N80").Value = -10542.875
$ws.Range("H83").Value = 6983.364
$ws.Range("I83").Value = 2814
$ws.Range("J83").Value = 8546.875
$ws.Range("K83").Value = 14070
$ws.Range("L83").Value = 42734.375
$ws.Range("M83").Value = -9078
$ws.Range("N83").Value = -52718.375
$ws.Range("H113").Value = 1820
$ws.Range("I113").Value = 1820
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1820
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 350
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 5853.6816
$ws.Range("I122").Value = 5410.3335
$ws.Range("J122").Value = 7848.75
$ws.Range("K122").Value = 16231.0005
$ws.Range("L122").Value = 23546.25
$ws.Range("M122").Value = -13781.0005
$ws.Range("N122").Value = -28446.25
$ws.Range("H126").Value = 6337
$ws.Range("I126").Value = 7255.5
$ws.Range("K126").Value = 21766.5
$ws.Range("M126").Value = -19296.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 999.5
$ws.Range("I7").Value = 999.5
$ws.Range("K7").Value = 999.5
$ws.Range("M7").Value = -887.5
$ws.Range("H40").Value = 76931640
$ws.Range("I40").Value = 111119080
$ws.Range("K40").Value = 111119080
$ws.Range("M40").Value = -111118944
$ws.Range("H46").Value = 1271
$ws.Range("I46").Value = 976
$ws.Range("J46").Value = 1566
$ws.Range("K46").Value = 976
$ws.Range("L46").Value = 1566
$ws.Range("M46").Value = -788
$ws.Range("N46").Value = -1942
$ws.Range("H55").Value = 946.3333
$ws.Range("I55").Value = 350.25
$ws.Range("J55").Value = 1313.1538
$ws.Range("K55").Value = 350.25
$ws.Range("L55").Value = 1313.1538
$ws.Range("M55").Value = -177.25
$ws.Range("N55").Value = -1659.1538
$ws.Range("H64").Value = 20075
$ws.Range("J64").Value = 20075
$ws.Range("L64").Value = 20075
$ws.Range("N64").Value = -20525
$ws.Range("H67").Value = 20075
$ws.Range("J67").Value = 20075
$ws.Range("L67").Value = 20075
$ws.Range("N67").Value = -21635
$ws.Range("H118").Value = 30000
$ws.Range("J118").Value = 30000
$ws.Range("L118").Value = 30000
$ws.Range("N118").Value = -33314
$ws.Range("H126").Value = 999.5
$ws.Range("I126").Value = 999.5
$ws.Range("K126").Value = 2998.5
$ws.Range("M126").Value = -528.5
$ws.Range("H132").Value = 2789.875
$ws.Range("I132").Value = 2789.875
$ws.Range("K132").Value = 8369.625
$ws.Range("M132").Value = -5839.625
$ws.Range("H136").Value = 3025.5557
$ws.Range("I136").Value = 3040.1428
$ws.Range("J136").Value = 2974.5
$ws.Range("K136").Value = 9120.428400000001
$ws.Range("L136").Value = 8923.5
$ws.Range("M136").Value = -6570.428400000001
$ws.Range("N136").Value = -14023.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 23839.572
$ws.Range("J74").Value = 22998.8
$ws.Range("L74").Value = 22998.8
$ws.Range("N74").Value = -24870.8
$ws.Range("H77").Value = 23839.572
$ws.Range("J77").Value = 22998.8
$ws.Range("L77").Value = 68996.39999999999
$ws.Range("N77").Value = -78356.39999999999
$ws.Range("H122").Value = 5948.8335
$ws.Range("I122").Value = 5948.8335
$ws.Range("K122").Value = 17846.5005
$ws.Range("M122").Value = -15396.5005
$ws.Range("H132").Value = 200007600
$ws.Range("I132").Value = 10000
$ws.Range("K132").Value = 30000
$ws.Range("M132").Value = -27470
$ws.Range("H136").Value = 9649.333000000001
$ws.Range("I136").Value = 11697.333
$ws.Range("K136").Value = 35091.999
$ws.Range("M136").Value = -32541.999
